$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Project"
$ws.Range("B2").Value = "Google"
$ws.Range("B5").Hyperlinks.Delete()
$ws.Range("B5").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("F8").Select() | Out-Null
